$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.370913
$ws.Range("H2").Value = 1.112739
$ws.Range("M2").Value = 15.50220733333333
$ws.Range("N2").Value = 46.506622
$ws.Range("O2").Value = 0.5994675913188158
$ws.Range("P2").Value = 0.5994675913188158
$ws.Range("Q2").Value = 5.749970228628666
$ws.Range("R2").Value = 51.74973205765799
$ws.Range("S2").Value = 0.5994675913188158
$ws.Range("T2").Value = 0.5994675913188158

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.370913
$ws.Range("H3").Value = 1.112739
$ws.Range("O3").Value = 0.04399860030713892
$ws.Range("P3").Value = 0.04399860030713892
$ws.Range("Q3").Value = 0.422025553226
$ws.Range("R3").Value = 3.798229979034
$ws.Range("S3").Value = 0.04399860030713892
$ws.Range("T3").Value = 0.04399860030713892

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.370913
$ws.Range("H4").Value = 1.112739
$ws.Range("M4").Value = 8.848210666666667
$ws.Range("N4").Value = 26.544632
$ws.Range("O4").Value = 0.3421587275782868
$ws.Range("P4").Value = 0.3421587275782868
$ws.Range("Q4").Value = 3.281916363005333
$ws.Range("R4").Value = 29.537247267048
$ws.Range("S4").Value = 0.3421587275782868
$ws.Range("T4").Value = 0.3421587275782868

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.370913
$ws.Range("H5").Value = 1.112739
$ws.Range("M5").Value = 0.371739
$ws.Range("N5").Value = 1.115217
$ws.Range("O5").Value = 0.01437508079575842
$ws.Range("P5").Value = 0.01437508079575841
$ws.Range("Q5").Value = 0.137882827707
$ws.Range("R5").Value = 1.240945449363
$ws.Range("S5").Value = 0.01437508079575842
$ws.Range("T5").Value = 0.01437508079575841
